{"js": "// The document contains a shell-command code block (notes on building\n// mod_wsgi from source). The real content edit is on the code line\n//   $ ./configure --with-python=python3\n// which becomes\n//   $ sudo ./configure --with-python=python3\n//\n// A handful of neighbouring code-block paragraphs also have runs that were\n// re-split/re-merged with no visible text change (the \"tar -xzf ...\" /\n// \"cd mod_wsgi-4.7.1/\" / \"make\" / \"make install\" lines). We reproduce those\n// run merges too, using narrow, paragraph-scoped search+replace so we never\n// touch more text than intended and never disturb paragraph formatting.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nfunction findParagraph(predicate, label) {\n  const match = paragraphs.items.find((p) => predicate(p.text));\n  if (!match) {\n    throw new Error(\"Could not locate paragraph: \" + label);\n  }\n  return match;\n}\n\nasync function replaceInParagraph(paragraph, find, replacement) {\n  const results = paragraph.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"$ sudo tar -xzf 4.7.1.tar.gz mod_wsgi-4.7.1/\"\nconst tarParagraph = findParagraph(\n  (t) => t.indexOf(\"tar -xzf\") !== -1,\n  \"tar -xzf ... command line\"\n);\nawait replaceInParagraph(tarParagraph, \" tar -\", \" tar -\");\nawait replaceInParagraph(\n  tarParagraph,\n  \" 4.7.1.tar.gz mod_wsgi-4.7.1/\",\n  \" 4.7.1.tar.gz mod_wsgi-4.7.1/\"\n);\n\n// \"$ cd mod_wsgi-4.7.1/\"\nconst cdParagraph = findParagraph(\n  (t) => t.indexOf(\"cd mod_wsgi-4.7.1/\") !== -1,\n  \"cd mod_wsgi-4.7.1/ command line\"\n);\nawait replaceInParagraph(cdParagraph, \"$ cd mod_wsgi-4.7.1/\", \"$ cd mod_wsgi-4.7.1/\");\n\n// \"$ ./configure --with-python=python3\" -> \"$ sudo ./configure --with-python=python3\"\nconst configureParagraph = findParagraph(\n  (t) => t.indexOf(\"./configure --with-python\") !== -1,\n  \"./configure --with-python=python3 command line\"\n);\nawait replaceInParagraph(configureParagraph, \"$ .\", \"$ sudo .\");\n\n// \"$ sudo make\" (not \"make install\")\nconst makeParagraph = findParagraph(\n  (t) => /\\$ sudo make$/.test(t.trim()),\n  \"$ sudo make command line\"\n);\nawait replaceInParagraph(makeParagraph, \" make\", \" make\");\n\n// \"$ sudo make install\"\nconst makeInstallParagraph = findParagraph(\n  (t) => t.indexOf(\"make install\") !== -1,\n  \"$ sudo make install command line\"\n);\nawait replaceInParagraph(makeInstallParagraph, \" make install\", \" make install\");\n", "ps1": "# The document contains a shell-command code block (notes on building\n# mod_wsgi from source). The real content edit is on the code line\n#   $ ./configure --with-python=python3\n# which becomes\n#   $ sudo ./configure --with-python=python3\n#\n# A handful of neighbouring code-block paragraphs also have runs that were\n# re-split/re-merged with no visible text change (the \"tar -xzf ...\" /\n# \"cd mod_wsgi-4.7.1/\" / \"make\" / \"make install\" lines). We reproduce those\n# run merges too, using narrow Find/Replace calls so we never touch more\n# text than intended and never disturb paragraph formatting.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\nfunction Replace-Text($searchText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $matchCase = $true\n    $ok = $find.Execute($searchText, $matchCase, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n    if (-not $ok) {\n        throw \"Find/Replace failed for '$searchText'\"\n    }\n}\n\n# \"$ sudo tar -xzf 4.7.1.tar.gz mod_wsgi-4.7.1/\"\nReplace-Text \" tar -\" \" tar -\"\nReplace-Text \" 4.7.1.tar.gz mod_wsgi-4.7.1/\" \" 4.7.1.tar.gz mod_wsgi-4.7.1/\"\n\n# \"$ cd mod_wsgi-4.7.1/\"\nReplace-Text \"`$ cd mod_wsgi-4.7.1/\" \"`$ cd mod_wsgi-4.7.1/\"\n\n# \"$ ./configure --with-python=python3\" -> \"$ sudo ./configure --with-python=python3\"\nReplace-Text \"`$ .\" \"`$ sudo .\"\n\n# \"$ sudo make\"\nReplace-Text \" make\" \" make\"\n\n# \"$ sudo make install\"\nReplace-Text \" make install\" \" make install\"\n"}
